$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1200
$ws.Range("I40").Value = 1000
$ws.Range("K40").Value = 1000
$ws.Range("M40").Value = -825

$ws.Range("H64").Value = 2755.5557
$ws.Range("I64").Value = 2740
$ws.Range("J64").Value = 2775
$ws.Range("K64").Value = 2740
$ws.Range("L64").Value = 2775
$ws.Range("M64").Value = -2492
$ws.Range("N64").Value = -3271

$ws.Range("H67").Value = 2755.5557
$ws.Range("I67").Value = 2740
$ws.Range("J67").Value = 2775
$ws.Range("K67").Value = 2740
$ws.Range("L67").Value = 2775
$ws.Range("M67").Value = -1882
$ws.Range("N67").Value = -4491

$ws.Range("H76").Value = 3145.7576
$ws.Range("I76").Value = 3096.6428
$ws.Range("J76").Value = 3420.8
$ws.Range("K76").Value = 3096.6428
$ws.Range("L76").Value = 3420.8
$ws.Range("M76").Value = -2781.6428
$ws.Range("N76").Value = -4050.8

$ws.Range("H79").Value = 3145.7576
$ws.Range("I79").Value = 3096.6428
$ws.Range("J79").Value = 3420.8
$ws.Range("K79").Value = 3096.6428
$ws.Range("L79").Value = 3420.8
$ws.Range("M79").Value = -2004.6428
$ws.Range("N79").Value = -5604.8

$ws.Range("H132").Value = 29417386
$ws.Range("I132").Value = 35720132
$ws.Range("J132").Value = 4567.5
$ws.Range("K132").Value = 107160396
$ws.Range("L132").Value = 13702.5
$ws.Range("M132").Value = -107157866
$ws.Range("N132").Value = -18762.5

$ws.Range("H137").Value = 3502.7886
$ws.Range("I137").Value = 3270.5
$ws.Range("K137").Value = 9811.5
$ws.Range("M137").Value = -7261.5

$ws.Range("H138").Value = 3733
$ws.Range("I138").Value = 1622
$ws.Range("J138").Value = 4355.918
$ws.Range("K138").Value = 4866
$ws.Range("L138").Value = 13067.754
$ws.Range("M138").Value = 274
$ws.Range("N138").Value = -23347.754

$ws.Range("H141").Value = 2452.76
$ws.Range("I141").Value = 2032.55
$ws.Range("J141").Value = 4133.6
$ws.Range("K141").Value = 6097.65
$ws.Range("L141").Value = 12400.8
$ws.Range("M141").Value = -917.6499999999996
$ws.Range("N141").Value = -22760.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14169.6
$ws.Range("I32").Value = 10780.857
$ws.Range("J32").Value = 18123.133
$ws.Range("K32").Value = 10780.857
$ws.Range("L32").Value = 18123.133
$ws.Range("M32").Value = -10493.857
$ws.Range("N32").Value = -18697.133

$ws.Range("H45").Value = 1023.7778
$ws.Range("I45").Value = 901.75
$ws.Range("J45").Value = 1121.4
$ws.Range("K45").Value = 901.75
$ws.Range("L45").Value = 1121.4
$ws.Range("M45").Value = -524.75
$ws.Range("N45").Value = -1875.4

$ws.Range("H61").Value = 3488.2222
$ws.Range("I61").Value = 2756.2856
$ws.Range("K61").Value = 2756.2856
$ws.Range("M61").Value = -2544.2856

$ws.Range("H122").Value = 3194.3333
$ws.Range("I122").Value = 1681.7273
$ws.Range("J122").Value = 5571.2856
$ws.Range("K122").Value = 5045.1819
$ws.Range("L122").Value = 16713.8568
$ws.Range("M122").Value = -2595.1819
$ws.Range("N122").Value = -21613.8568

$ws.Range("H136").Value = 3488.2222
$ws.Range("I136").Value = 2756.2856
$ws.Range("K136").Value = 8268.856800000001
$ws.Range("M136").Value = -5718.856800000001

$ws.Range("H140").Value = 65143
$ws.Range("J140").Value = 65143
$ws.Range("L140").Value = 65143
$ws.Range("N140").Value = -75503

$ws.Range("H141").Value = 88334.836
$ws.Range("J141").Value = 88334.836
$ws.Range("L141").Value = 88334.836
$ws.Range("N141").Value = -98694.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2241.742
$ws.Range("I86").Value = 1740.8334
$ws.Range("J86").Value = 2558.1052
$ws.Range("K86").Value = 1740.8334
$ws.Range("L86").Value = 2558.1052
$ws.Range("M86").Value = -617.8334
$ws.Range("N86").Value = -4804.1052

$ws.Range("H89").Value = 2241.742
$ws.Range("I89").Value = 1740.8334
$ws.Range("J89").Value = 2558.1052
$ws.Range("K89").Value = 8704.166999999999
$ws.Range("L89").Value = 12790.526
$ws.Range("M89").Value = -3088.166999999999
$ws.Range("N89").Value = -24022.526

$ws.Range("H140").Value = 48025.668
$ws.Range("J140").Value = 48025.668
$ws.Range("L140").Value = 48025.668
$ws.Range("N140").Value = -58385.668

$ws.Range("H141").Value = 32000
$ws.Range("J141").Value = 32000
$ws.Range("L141").Value = 32000
$ws.Range("N141").Value = -42360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1323.6666
$ws.Range("I16").Value = 1129.4117
$ws.Range("J16").Value = 1653.9
$ws.Range("K16").Value = 1129.4117
$ws.Range("L16").Value = 1653.9
$ws.Range("M16").Value = -842.4117000000001
$ws.Range("N16").Value = -2227.9

$ws.Range("H62").Value = 35719612
$ws.Range("I62").Value = 55560544
$ws.Range("J62").Value = 5933.8
$ws.Range("K62").Value = 55560544
$ws.Range("L62").Value = 5933.8
$ws.Range("M62").Value = -55559920
$ws.Range("N62").Value = -7181.8

$ws.Range("H65").Value = 35719612
$ws.Range("I65").Value = 55560544
$ws.Range("J65").Value = 5933.8
$ws.Range("K65").Value = 277802720
$ws.Range("L65").Value = 29669
$ws.Range("M65").Value = -277799600
$ws.Range("N65").Value = -35909

$ws.Range("H107").Value = 439.35294
$ws.Range("I107").Value = 385.5625
$ws.Range("K107").Value = 385.5625
$ws.Range("M107").Value = 1534.4375

$ws.Range("H113").Value = 1323.6666
$ws.Range("I113").Value = 1129.4117
$ws.Range("J113").Value = 1653.9
$ws.Range("K113").Value = 1129.4117
$ws.Range("L113").Value = 1653.9
$ws.Range("M113").Value = 1040.5883
$ws.Range("N113").Value = -5993.9

$ws.Range("H134").Value = 5507.5757
$ws.Range("I134").Value = 5590.4443
$ws.Range("K134").Value = 16771.3329
$ws.Range("M134").Value = -14236.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7940668.5
$ws.Range("I131").Value = 26327230
$ws.Range("J131").Value = 1016.4773
$ws.Range("K131").Value = 78981690
$ws.Range("L131").Value = 3049.4319
$ws.Range("M131").Value = -78976650
$ws.Range("N131").Value = -13129.4319

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 31707
$ws.Range("J64").Value = 31707
$ws.Range("L64").Value = 31707
$ws.Range("N64").Value = -32203

$ws.Range("H67").Value = 31707
$ws.Range("J67").Value = 31707
$ws.Range("L67").Value = 31707
$ws.Range("N67").Value = -33423

$ws.Range("H70").Value = 6167.759
$ws.Range("J70").Value = 7343.4
$ws.Range("L70").Value = 7343.4
$ws.Range("N70").Value = -7883.4

$ws.Range("H73").Value = 6167.759
$ws.Range("J73").Value = 7343.4
$ws.Range("L73").Value = 7343.4
$ws.Range("N73").Value = -9215.4

$ws.Range("H80").Value = 50002300
$ws.Range("I80").Value = 83335170
$ws.Range("K80").Value = 83335170
$ws.Range("M80").Value = -83334172

$ws.Range("H83").Value = 50002300
$ws.Range("I83").Value = 83335170
$ws.Range("K83").Value = 416675850
$ws.Range("M83").Value = -416670858

$ws.Range("H122").Value = 2628.8057
$ws.Range("I122").Value = 1964.091
$ws.Range("K122").Value = 5892.272999999999
$ws.Range("M122").Value = -3442.272999999999

$ws.Range("H126").Value = 4099.298
$ws.Range("I126").Value = 2975.9023
$ws.Range("J126").Value = 5170.442
$ws.Range("K126").Value = 8927.706900000001
$ws.Range("L126").Value = 15511.326
$ws.Range("M126").Value = -6457.706900000001
$ws.Range("N126").Value = -20451.326

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3514.5144
$ws.Range("I122").Value = 2454.2917
$ws.Range("K122").Value = 7362.875100000001
$ws.Range("M122").Value = -4912.875100000001
